$wb = $excel.ActiveWorkbook

$oldText = "ARIMA model for the time series (The second mortality index of the reduced Plat model with two factors)"
$newText = "ARIMA model for the time series (the second mortality index of the reduced Plat model with two factors)"

foreach ($ws in $wb.Worksheets) {
    $cell = $ws.Range("A1")
    if ($cell.Value2 -eq $oldText) {
        $cell.Value = $newText
    }
}
